# Insert two new price-record rows (392 and 393) into the Brócoli sheet,
# pushing the previously-existing rows 392:493 down to 394:495.
#
# Row 392 (new):  21-12-2021, Primera, Volumen 1908, Precio min/max 600/600, Precio promedio 600
# Row 393 (new):  21-12-2021, Segunda, Volumen 1950, Precio min/max 500/500, Precio promedio 500

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 392 (shifts old rows 392:493 -> 394:495).
$ws.Range("A392:A393").EntireRow.Insert()

# Seed the two new rows with the "constant" columns (Mercado, Región, Codreg,
# Categoría, Variedad, Unidad de comercialización, Origen, Kg o Unidades,
# Clasificación, etc.) by copying them from the rows that now sit directly
# below (which hold what used to be rows 392/393 before the shift, and carry
# identical constant-column values for this market/category).
$ws.Range("A392:R392").Value = $ws.Range("A394:R394").Value()
$ws.Range("A393:R393").Value = $ws.Range("A395:R395").Value()

# Overwrite the record-specific fields for the first new row (392).
$ws.Range("D392").Value = "2021-12-21"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 1908
$ws.Range("K392").Value = 600
$ws.Range("L392").Value = 600
$ws.Range("M392").Value = 600
$ws.Range("P392").Value = 600

# Overwrite the record-specific fields for the second new row (393).
$ws.Range("D393").Value = "2021-12-21"
$ws.Range("I393").Value = "Segunda"
$ws.Range("J393").Value = 1950
$ws.Range("K393").Value = 500
$ws.Range("L393").Value = 500
$ws.Range("M393").Value = 500
$ws.Range("P393").Value = 500
